$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Log_Muestras")

$timestamps = @{
    2  = "2025-11-02T02:03:23.683903"
    3  = "2025-11-02T02:03:23.683903"
    4  = "2025-11-02T02:03:23.683903"
    5  = "2025-11-02T02:03:23.684529"
    6  = "2025-11-02T02:03:23.684529"
    7  = "2025-11-02T02:03:23.684529"
    8  = "2025-11-02T02:03:23.685075"
    9  = "2025-11-02T02:03:23.685075"
    10 = "2025-11-02T02:03:23.686190"
    11 = "2025-11-02T02:03:23.686190"
    12 = "2025-11-02T02:03:23.686190"
    13 = "2025-11-02T02:03:23.686190"
    14 = "2025-11-02T02:03:23.686190"
    15 = "2025-11-02T02:03:23.686190"
    16 = "2025-11-02T02:03:23.686190"
    17 = "2025-11-02T02:03:23.686190"
    18 = "2025-11-02T02:03:23.686190"
    19 = "2025-11-02T02:03:23.686190"
    20 = "2025-11-02T02:03:23.686190"
    21 = "2025-11-02T02:03:23.686190"
    22 = "2025-11-02T02:03:23.689327"
    23 = "2025-11-02T02:03:23.689753"
    24 = "2025-11-02T02:03:23.689753"
    25 = "2025-11-02T02:03:23.690290"
    26 = "2025-11-02T02:03:23.690290"
    27 = "2025-11-02T02:03:23.690290"
    28 = "2025-11-02T02:03:23.690826"
    29 = "2025-11-02T02:03:23.690826"
}

foreach ($row in $timestamps.Keys) {
    $ws.Cells.Item($row, 26).Value = $timestamps[$row]
}
